$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 23
$ws1.Range("F4").Value = 239
$ws1.Range("F6").Value = 1152
$ws1.Range("F7").Value = 932
$ws1.Range("F11").Value = 894
$ws1.Range("F12").Value = 321
$ws1.Range("F15").Value = 1375
$ws1.Range("F17").Value = 1273
$ws1.Range("F18").Value = 2935
$ws1.Range("F19").Value = 261
$ws1.Range("F20").Value = 1557
$ws1.Range("F21").Value = 1311
$ws1.Range("F23").Value = 216
$ws1.Range("F26").Value = 1071
$ws1.Range("F28").Value = 3308
$ws1.Range("F31").Value = 1468

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 37
$ws2.Range("F10").Value = 7

# Sheet: 全部类型 (All types) - mirrors the combined data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 23
$ws4.Range("F7").Value = 239
$ws4.Range("F10").Value = 1152
$ws4.Range("F11").Value = 932
$ws4.Range("F19").Value = 37
$ws4.Range("F21").Value = 7
$ws4.Range("F23").Value = 894
$ws4.Range("F24").Value = 321
$ws4.Range("F27").Value = 1375
$ws4.Range("F29").Value = 1273
$ws4.Range("F30").Value = 2935
$ws4.Range("F31").Value = 261
$ws4.Range("F32").Value = 1557
$ws4.Range("F33").Value = 1311
$ws4.Range("F35").Value = 216
$ws4.Range("F40").Value = 1071
$ws4.Range("F42").Value = 3308
$ws4.Range("F45").Value = 1468
